# Insert two new data rows (weekly price observations) into the Albahaca
# sheet, right after the existing row 519, shifting all subsequent rows
# down by two (old row N -> new row N+2 for N = 520..630).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 520 (Excel shifts 520:630 down to 522:632,
# picking up the formatting of the row above for the new rows).
$ws.Rows.Item(520).Resize(2).Insert()

# --- New row 520 -----------------------------------------------------
$ws.Range("A520").Value = 9
$ws.Range("B520").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C520").Value = "Metropolitana"
$ws.Range("D520").Value = 45211
$ws.Range("E520").Value = 13
$ws.Range("F520").Value = 100112052
$ws.Range("G520").Value = "Albahaca"
$ws.Range("H520").Value = "Sin especificar"
$ws.Range("I520").Value = "Primera"
$ws.Range("J520").Value = 70
$ws.Range("K520").Value = 5000
$ws.Range("L520").Value = 5000
$ws.Range("M520").Value = 5000
$ws.Range("N520").Value = "`$/docena de matas"
$ws.Range("O520").Value = "Provincia de Chacabuco"
$ws.Range("P520").Value = 833
$ws.Range("Q520").Value = 6
$ws.Range("R520").Value = "Hortaliza"

# --- New row 521 -----------------------------------------------------
$ws.Range("A521").Value = 9
$ws.Range("B521").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C521").Value = "Metropolitana"
$ws.Range("D521").Value = 45211
$ws.Range("E521").Value = 13
$ws.Range("F521").Value = 100112052
$ws.Range("G521").Value = "Albahaca"
$ws.Range("H521").Value = "Sin especificar"
$ws.Range("I521").Value = "Primera"
$ws.Range("J521").Value = 430
$ws.Range("K521").Value = 5000
$ws.Range("L521").Value = 5500
$ws.Range("M521").Value = 5250
$ws.Range("N521").Value = "`$/paquete"
$ws.Range("O521").Value = "Región de Arica y Parinacota"
$ws.Range("P521").Value = 5250
$ws.Range("Q521").Value = 1
$ws.Range("R521").Value = "Hortaliza"

Write-Output "Inserted 2 rows and populated D520:R521"
